$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44708
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12571
$ws.Range("Q2").Value = '$/caja 12 kilos empedrada'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1048
$ws.Range("T2").Value = 12

$ws.Range("D3").Value = 44714
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 806
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 45106
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17625
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región del Maule'
$ws.Range("S4").Value = 979
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 45114
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Región del Maule'
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 45092
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 140
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18429
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 1024
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44742
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 806
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44330
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 861
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 45084
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 17500
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región del Maule'
$ws.Range("S9").Value = 972
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44719
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14400
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región del Maule'
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 45090
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 150
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 17533
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Región del Maule'
$ws.Range("S11").Value = 974
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 45090
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 130
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14462
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región del Maule'
$ws.Range("S12").Value = 803
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 45091
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 220
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 19000
$ws.Range("P13").Value = 18455
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 1025
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = 45091
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 833
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44707
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 12500
$ws.Range("Q15").Value = '$/caja 12 kilos empedrada'
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 1042
$ws.Range("T15").Value = 12

$ws.Range("D16").Value = 45097
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región del Maule'
$ws.Range("S16").Value = 1056
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 45100
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 18000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 18000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1000
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 44334
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 11000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 11500
$ws.Range("Q18").Value = '$/caja 12 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 11500
$ws.Range("T18").Value = 1

$ws.Range("D19").Value = 45077
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 140
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 12857
$ws.Range("Q19").Value = '$/caja 12 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 12857
$ws.Range("T19").Value = 1

$ws.Range("D20").Value = 45077
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 11000
$ws.Range("P20").Value = 11000
$ws.Range("Q20").Value = '$/caja 12 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 11000
$ws.Range("T20").Value = 1

$ws.Range("D21").Value = 45093
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 140
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 17429
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 968
$ws.Range("T21").Value = 18
